$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Config_4")
$ws.Cells.Item(2, 2).Value = 2499
$ws.Cells.Item(2, 3).Value = 751
$ws.Cells.Item(2, 4).Value = 85.269032197
$ws.Cells.Item(4, 2).Value = 2494
$ws.Cells.Item(4, 3).Value = 877
$ws.Cells.Item(4, 4).Value = 100.477497363
$ws.Cells.Item(5, 2).Value = 3240
$ws.Cells.Item(5, 3).Value = 694
$ws.Cells.Item(5, 4).Value = 83.184140658
$ws.Cells.Item(6, 2).Value = 3414
$ws.Cells.Item(6, 3).Value = 901
$ws.Cells.Item(6, 4).Value = 110.669556724
$ws.Cells.Item(7, 2).Value = 3060
$ws.Cells.Item(7, 3).Value = 611
$ws.Cells.Item(7, 4).Value = 73.215914181
$ws.Cells.Item(8, 2).Value = 3403
$ws.Cells.Item(8, 3).Value = 855
$ws.Cells.Item(8, 4).Value = 105.278175504
$ws.Cells.Item(9, 2).Value = 3575
$ws.Cells.Item(9, 3).Value = 902
$ws.Cells.Item(9, 4).Value = 111.028542434
$ws.Cells.Item(10, 2).Value = 3238
$ws.Cells.Item(10, 3).Value = 920
$ws.Cells.Item(10, 4).Value = 108.944551003
$ws.Cells.Item(11, 2).Value = 3510
$ws.Cells.Item(11, 3).Value = 1389
$ws.Cells.Item(11, 4).Value = 166.736322464

$ws = $wb.Worksheets.Item("Config_6")
$ws.Cells.Item(2, 2).Value = 2499
$ws.Cells.Item(2, 3).Value = 774
$ws.Cells.Item(2, 4).Value = 88.030289811
$ws.Cells.Item(4, 2).Value = 2494
$ws.Cells.Item(4, 3).Value = 877
$ws.Cells.Item(4, 4).Value = 103.173261302
$ws.Cells.Item(5, 2).Value = 3240
$ws.Cells.Item(5, 3).Value = 694
$ws.Cells.Item(5, 4).Value = 83.324901231
$ws.Cells.Item(6, 2).Value = 3414
$ws.Cells.Item(6, 3).Value = 901
$ws.Cells.Item(6, 4).Value = 110.510246917
$ws.Cells.Item(7, 2).Value = 3060
$ws.Cells.Item(7, 3).Value = 611
$ws.Cells.Item(7, 4).Value = 73.255368101
$ws.Cells.Item(8, 2).Value = 3403
$ws.Cells.Item(8, 3).Value = 855
$ws.Cells.Item(8, 4).Value = 105.559537059
$ws.Cells.Item(9, 2).Value = 3575
$ws.Cells.Item(9, 3).Value = 902
$ws.Cells.Item(9, 4).Value = 111.184486483
$ws.Cells.Item(10, 2).Value = 3238
$ws.Cells.Item(10, 3).Value = 1979
$ws.Cells.Item(10, 4).Value = 229.925217957
$ws.Cells.Item(11, 2).Value = 3510
$ws.Cells.Item(11, 3).Value = 1389
$ws.Cells.Item(11, 4).Value = 166.54630782

$ws = $wb.Worksheets.Item("Config_11")
$ws.Cells.Item(2, 2).Value = 3202
$ws.Cells.Item(2, 3).Value = 3304
$ws.Cells.Item(2, 4).Value = 465.005538699
$ws.Cells.Item(3, 2).Value = 3254
$ws.Cells.Item(3, 3).Value = 3824
$ws.Cells.Item(3, 4).Value = 550.518663999
$ws.Cells.Item(4, 2).Value = 3483
$ws.Cells.Item(4, 3).Value = 3128
$ws.Cells.Item(4, 4).Value = 454.693188016
$ws.Cells.Item(5, 2).Value = 3873
$ws.Cells.Item(5, 3).Value = 3073
$ws.Cells.Item(5, 4).Value = 459.408027875
$ws.Cells.Item(6, 3).Value = 2977
$ws.Cells.Item(6, 4).Value = 406.75001835
$ws.Cells.Item(7, 2).Value = 4007
$ws.Cells.Item(7, 3).Value = 2580
$ws.Cells.Item(7, 4).Value = 390.636497441
$ws.Cells.Item(8, 2).Value = 2380
$ws.Cells.Item(8, 3).Value = 2926
$ws.Cells.Item(8, 4).Value = 445.194711279
$ws.Cells.Item(9, 2).Value = 3833
$ws.Cells.Item(9, 3).Value = 2598
$ws.Cells.Item(9, 4).Value = 387.089492451
$ws.Cells.Item(10, 2).Value = 3438
$ws.Cells.Item(10, 3).Value = 3020
$ws.Cells.Item(10, 4).Value = 462.233871756
$ws.Cells.Item(11, 2).Value = 2580
$ws.Cells.Item(11, 3).Value = 3763
$ws.Cells.Item(11, 4).Value = 587.719382232
$ws.Cells.Item(12, 2).Value = 2640
$ws.Cells.Item(12, 3).Value = 3972
$ws.Cells.Item(12, 4).Value = 522.293294571
$ws.Cells.Item(13, 2).Value = 3456
$ws.Cells.Item(13, 3).Value = 3758
$ws.Cells.Item(13, 4).Value = 529.964491217
$ws.Cells.Item(14, 2).Value = 3389
$ws.Cells.Item(14, 3).Value = 3535
$ws.Cells.Item(14, 4).Value = 517.741827205
$ws.Cells.Item(15, 2).Value = 3099
$ws.Cells.Item(15, 3).Value = 3088
$ws.Cells.Item(15, 4).Value = 449.149723489
$ws.Cells.Item(16, 2).Value = 3235
$ws.Cells.Item(16, 3).Value = 3131
$ws.Cells.Item(16, 4).Value = 476.28905885
$ws.Cells.Item(17, 2).Value = 3926
$ws.Cells.Item(17, 3).Value = 2698
$ws.Cells.Item(17, 4).Value = 369.410790214
$ws.Cells.Item(18, 2).Value = 3241
$ws.Cells.Item(18, 3).Value = 3203
$ws.Cells.Item(18, 4).Value = 470.99681426
$ws.Cells.Item(19, 2).Value = 3155
$ws.Cells.Item(19, 3).Value = 3101
$ws.Cells.Item(19, 4).Value = 449.032495451
$ws.Cells.Item(20, 2).Value = 3147
$ws.Cells.Item(20, 3).Value = 3980
$ws.Cells.Item(20, 4).Value = 526.229015058

$ws = $wb.Worksheets.Item("Summary_n_tokens_prompt")
$ws.Cells.Item(4, 3).Value = 3089.9
$ws.Cells.Item(4, 4).Value = 3239
$ws.Cells.Item(4, 5).Value = 441.4303644592958
$ws.Cells.Item(4, 7).Value = 3575
$ws.Cells.Item(4, 8).Value = 2639.25
$ws.Cells.Item(4, 9).Value = 3411.25
$ws.Cells.Item(4, 11).Value = 30899
$ws.Cells.Item(6, 3).Value = 3089.9
$ws.Cells.Item(6, 4).Value = 3239
$ws.Cells.Item(6, 5).Value = 441.4303644592958
$ws.Cells.Item(6, 7).Value = 3575
$ws.Cells.Item(6, 8).Value = 2639.25
$ws.Cells.Item(6, 9).Value = 3411.25
$ws.Cells.Item(6, 11).Value = 30899
$ws.Cells.Item(9, 3).Value = 3317.631578947368
$ws.Cells.Item(9, 4).Value = 3254
$ws.Cells.Item(9, 5).Value = 450.0182972238543
$ws.Cells.Item(9, 6).Value = 2380
$ws.Cells.Item(9, 7).Value = 4007
$ws.Cells.Item(9, 8).Value = 3151
$ws.Cells.Item(9, 9).Value = 3590
$ws.Cells.Item(9, 11).Value = 63035

$ws = $wb.Worksheets.Item("Summary_n_tokens_response")
$ws.Cells.Item(4, 3).Value = 874.3
$ws.Cells.Item(4, 4).Value = 866
$ws.Cells.Item(4, 5).Value = 207.3328242223117
$ws.Cells.Item(4, 6).Value = 611
$ws.Cells.Item(4, 7).Value = 1389
$ws.Cells.Item(4, 8).Value = 774
$ws.Cells.Item(4, 9).Value = 901.75
$ws.Cells.Item(4, 11).Value = 8743
$ws.Cells.Item(6, 3).Value = 983
$ws.Cells.Item(6, 4).Value = 866
$ws.Cells.Item(6, 5).Value = 405.7322598298866
$ws.Cells.Item(6, 6).Value = 611
$ws.Cells.Item(6, 7).Value = 1979
$ws.Cells.Item(6, 8).Value = 792.5
$ws.Cells.Item(6, 9).Value = 901.75
$ws.Cells.Item(6, 11).Value = 9830
$ws.Cells.Item(9, 3).Value = 3245.21052631579
$ws.Cells.Item(9, 4).Value = 3128
$ws.Cells.Item(9, 5).Value = 441.634537063719
$ws.Cells.Item(9, 6).Value = 2580
$ws.Cells.Item(9, 7).Value = 3980
$ws.Cells.Item(9, 8).Value = 2998.5
$ws.Cells.Item(9, 9).Value = 3646.5
$ws.Cells.Item(9, 11).Value = 61659

$ws = $wb.Worksheets.Item("Summary_time")
$ws.Cells.Item(4, 3).Value = 117.4049190889
$ws.Cells.Item(4, 4).Value = 107.1113632535
$ws.Cells.Item(4, 5).Value = 46.77313482351462
$ws.Cells.Item(4, 6).Value = 73.215914181
$ws.Cells.Item(4, 7).Value = 229.245458361
$ws.Cells.Item(4, 8).Value = 89.0711484885
$ws.Cells.Item(4, 9).Value = 110.9387960065
$ws.Cells.Item(4, 11).Value = 1174.049190889
$ws.Cells.Item(6, 3).Value = 122.5592382633
$ws.Cells.Item(6, 4).Value = 108.034891988
$ws.Cells.Item(6, 5).Value = 47.79338737685489
$ws.Cells.Item(6, 6).Value = 73.255368101
$ws.Cells.Item(6, 7).Value = 229.925217957
$ws.Cells.Item(6, 8).Value = 91.81603268375
$ws.Cells.Item(6, 9).Value = 143.35819608475
$ws.Cells.Item(6, 11).Value = 1225.592382633
$ws.Cells.Item(9, 3).Value = 469.4924685480526
$ws.Cells.Item(9, 4).Value = 462.233871756
$ws.Cells.Item(9, 5).Value = 58.38614679775192
$ws.Cells.Item(9, 6).Value = 369.410790214
$ws.Cells.Item(9, 7).Value = 587.719382232
$ws.Cells.Item(9, 8).Value = 447.113603365
$ws.Cells.Item(9, 9).Value = 520.017560888
$ws.Cells.Item(9, 11).Value = 8920.356902412999

$ws = $wb.Worksheets.Item("Raw_Data")
$ws.Cells.Item(25, 6).Value = 2499
$ws.Cells.Item(25, 7).Value = 751
$ws.Cells.Item(25, 8).Value = 85.269032197
$ws.Cells.Item(27, 6).Value = 2494
$ws.Cells.Item(27, 7).Value = 877
$ws.Cells.Item(27, 8).Value = 100.477497363
$ws.Cells.Item(28, 6).Value = 3240
$ws.Cells.Item(28, 7).Value = 694
$ws.Cells.Item(28, 8).Value = 83.184140658
$ws.Cells.Item(29, 6).Value = 3414
$ws.Cells.Item(29, 7).Value = 901
$ws.Cells.Item(29, 8).Value = 110.669556724
$ws.Cells.Item(30, 6).Value = 3060
$ws.Cells.Item(30, 7).Value = 611
$ws.Cells.Item(30, 8).Value = 73.215914181
$ws.Cells.Item(31, 6).Value = 3403
$ws.Cells.Item(31, 7).Value = 855
$ws.Cells.Item(31, 8).Value = 105.278175504
$ws.Cells.Item(32, 6).Value = 3575
$ws.Cells.Item(32, 7).Value = 902
$ws.Cells.Item(32, 8).Value = 111.028542434
$ws.Cells.Item(33, 6).Value = 3238
$ws.Cells.Item(33, 7).Value = 920
$ws.Cells.Item(33, 8).Value = 108.944551003
$ws.Cells.Item(34, 6).Value = 3510
$ws.Cells.Item(34, 7).Value = 1389
$ws.Cells.Item(34, 8).Value = 166.736322464
$ws.Cells.Item(58, 6).Value = 2499
$ws.Cells.Item(58, 7).Value = 774
$ws.Cells.Item(58, 8).Value = 88.030289811
$ws.Cells.Item(60, 6).Value = 2494
$ws.Cells.Item(60, 7).Value = 877
$ws.Cells.Item(60, 8).Value = 103.173261302
$ws.Cells.Item(61, 6).Value = 3240
$ws.Cells.Item(61, 7).Value = 694
$ws.Cells.Item(61, 8).Value = 83.324901231
$ws.Cells.Item(62, 6).Value = 3414
$ws.Cells.Item(62, 7).Value = 901
$ws.Cells.Item(62, 8).Value = 110.510246917
$ws.Cells.Item(63, 6).Value = 3060
$ws.Cells.Item(63, 7).Value = 611
$ws.Cells.Item(63, 8).Value = 73.255368101
$ws.Cells.Item(64, 6).Value = 3403
$ws.Cells.Item(64, 7).Value = 855
$ws.Cells.Item(64, 8).Value = 105.559537059
$ws.Cells.Item(65, 6).Value = 3575
$ws.Cells.Item(65, 7).Value = 902
$ws.Cells.Item(65, 8).Value = 111.184486483
$ws.Cells.Item(66, 6).Value = 3238
$ws.Cells.Item(66, 7).Value = 1979
$ws.Cells.Item(66, 8).Value = 229.925217957
$ws.Cells.Item(67, 6).Value = 3510
$ws.Cells.Item(67, 7).Value = 1389
$ws.Cells.Item(67, 8).Value = 166.54630782
$ws.Cells.Item(110, 6).Value = 3202
$ws.Cells.Item(110, 7).Value = 3304
$ws.Cells.Item(110, 8).Value = 465.005538699
$ws.Cells.Item(111, 6).Value = 3254
$ws.Cells.Item(111, 7).Value = 3824
$ws.Cells.Item(111, 8).Value = 550.518663999
$ws.Cells.Item(112, 6).Value = 3483
$ws.Cells.Item(112, 7).Value = 3128
$ws.Cells.Item(112, 8).Value = 454.693188016
$ws.Cells.Item(113, 6).Value = 3873
$ws.Cells.Item(113, 7).Value = 3073
$ws.Cells.Item(113, 8).Value = 459.408027875
$ws.Cells.Item(114, 7).Value = 2977
$ws.Cells.Item(114, 8).Value = 406.75001835
$ws.Cells.Item(115, 6).Value = 4007
$ws.Cells.Item(115, 7).Value = 2580
$ws.Cells.Item(115, 8).Value = 390.636497441
$ws.Cells.Item(116, 6).Value = 2380
$ws.Cells.Item(116, 7).Value = 2926
$ws.Cells.Item(116, 8).Value = 445.194711279
$ws.Cells.Item(117, 6).Value = 3833
$ws.Cells.Item(117, 7).Value = 2598
$ws.Cells.Item(117, 8).Value = 387.089492451
$ws.Cells.Item(118, 6).Value = 3438
$ws.Cells.Item(118, 7).Value = 3020
$ws.Cells.Item(118, 8).Value = 462.233871756
$ws.Cells.Item(119, 6).Value = 2580
$ws.Cells.Item(119, 7).Value = 3763
$ws.Cells.Item(119, 8).Value = 587.719382232
$ws.Cells.Item(120, 6).Value = 2640
$ws.Cells.Item(120, 7).Value = 3972
$ws.Cells.Item(120, 8).Value = 522.293294571
$ws.Cells.Item(121, 6).Value = 3456
$ws.Cells.Item(121, 7).Value = 3758
$ws.Cells.Item(121, 8).Value = 529.964491217
$ws.Cells.Item(122, 6).Value = 3389
$ws.Cells.Item(122, 7).Value = 3535
$ws.Cells.Item(122, 8).Value = 517.741827205
$ws.Cells.Item(123, 6).Value = 3099
$ws.Cells.Item(123, 7).Value = 3088
$ws.Cells.Item(123, 8).Value = 449.149723489
$ws.Cells.Item(124, 6).Value = 3235
$ws.Cells.Item(124, 7).Value = 3131
$ws.Cells.Item(124, 8).Value = 476.28905885
$ws.Cells.Item(125, 6).Value = 3926
$ws.Cells.Item(125, 7).Value = 2698
$ws.Cells.Item(125, 8).Value = 369.410790214
$ws.Cells.Item(126, 6).Value = 3241
$ws.Cells.Item(126, 7).Value = 3203
$ws.Cells.Item(126, 8).Value = 470.99681426
$ws.Cells.Item(127, 6).Value = 3155
$ws.Cells.Item(127, 7).Value = 3101
$ws.Cells.Item(127, 8).Value = 449.032495451
$ws.Cells.Item(128, 6).Value = 3147
$ws.Cells.Item(128, 7).Value = 3980
$ws.Cells.Item(128, 8).Value = 526.229015058
